# Last update with TARSO model
# Insert a new row above the existing "AR(4)" row (row 4) on the summary
# sheet and populate it with the new "Persistent" model results. This
# pushes every subsequent row down by one, matching the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4 (shifts rows 4.. down to 5..)
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the Persistent model's results.
# Column B is written before column A so new shared strings are
# appended in the same order as the source workbook.
$ws.Range("B4").Value = "past p 1h"
$ws.Range("A4").Value = "Persistent"
$ws.Range("C4").Value = "±3.25"
$ws.Range("D4").Value = "±4.88"
$ws.Range("E4").Value = "±6.03"

# Match the final selection left behind in the saved workbook.
$ws.Range("H2").Select()
